$d = $word.ActiveDocument
$nbsp = [char]0x00A0

# --- Paragraph 2: "Taxi by ... via Poly Pizza (<hyperlink>)" ---
# Collapse the three runs (plain text + hyperlinked URL + closing paren)
# into a single plain-text run, matching the target OOXML. The original
# run text uses non-breaking spaces in a few spots, so reproduce those
# verbatim instead of plain ASCII spaces.
$p2 = $d.Paragraphs.Item(2).Range
$p2.MoveEnd(1, -1)
$p2.Delete()
$p2b = $d.Paragraphs.Item(2).Range
$p2b.InsertAfter("Taxi${nbsp}by${nbsp}Poly by Google${nbsp}[CC-BY] (https://creativecommons.org/licenses/by/3.0/) via Poly Pizza (https://poly.pizza/m/2Me_E4PMM5J)")

# --- Remove the "BUS credits" paragraph and the Montreal Bus credits
# paragraph that followed it entirely. ---
$p3 = $d.Paragraphs.Item(3).Range
$p3.Delete()
$p3b = $d.Paragraphs.Item(3).Range
$p3b.Delete()

# --- Remove the now-unused character styles. ---
$sUnresolved = $d.Styles.Item("Unresolved Mention")
$sUnresolved.Delete()
$sHyperlink = $d.Styles.Item("Hyperlink")
$sHyperlink.Delete()
